# Weekly price update: a new week's worth of data (two rows: "Primera" and
# "Segunda" quality grades) is inserted into the logged series right after
# the first block of rows, pushing every subsequent row down by two and
# growing the used range from A1:R435 to A1:R437.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at 380:381, shifting rows 380:435 down to 382:437.
$ws.Rows("380:381").Insert(-4121)

# Populate the newly inserted row 380 ("Primera" quality).
$ws.Cells.Item(380, 1).Value = 8
$ws.Cells.Item(380, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(380, 3).Value = "Coquimbo"
$ws.Cells.Item(380, 4).Value = 44748
$ws.Cells.Item(380, 5).Value = 4
$ws.Cells.Item(380, 6).Value = 100112009
$ws.Cells.Item(380, 7).Value = "Acelga"
$ws.Cells.Item(380, 8).Value = "Sin especificar"
$ws.Cells.Item(380, 9).Value = "Primera"
$ws.Cells.Item(380, 10).Value = 2560
$ws.Cells.Item(380, 11).Value = 600
$ws.Cells.Item(380, 12).Value = 700
$ws.Cells.Item(380, 13).Value = 650
$ws.Cells.Item(380, 14).Value = "`$/atado 1,5 a 2 kilos"
$ws.Cells.Item(380, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(380, 16).Value = 325
$ws.Cells.Item(380, 17).Value = 2
$ws.Cells.Item(380, 18).Value = "Hortaliza"

# Populate the newly inserted row 381 ("Segunda" quality).
$ws.Cells.Item(381, 1).Value = 8
$ws.Cells.Item(381, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(381, 3).Value = "Coquimbo"
$ws.Cells.Item(381, 4).Value = 44748
$ws.Cells.Item(381, 5).Value = 4
$ws.Cells.Item(381, 6).Value = 100112009
$ws.Cells.Item(381, 7).Value = "Acelga"
$ws.Cells.Item(381, 8).Value = "Sin especificar"
$ws.Cells.Item(381, 9).Value = "Segunda"
$ws.Cells.Item(381, 10).Value = 1400
$ws.Cells.Item(381, 11).Value = 500
$ws.Cells.Item(381, 12).Value = 550
$ws.Cells.Item(381, 13).Value = 525
$ws.Cells.Item(381, 14).Value = "`$/atado 1,5 a 2 kilos"
$ws.Cells.Item(381, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(381, 16).Value = 262
$ws.Cells.Item(381, 17).Value = 2
$ws.Cells.Item(381, 18).Value = "Hortaliza"

# Make sure the date cells keep the date display format used by the rest of
# column D (style carried over from the Insert, but set explicitly to be safe).
$ws.Range("D380").NumberFormat = $ws.Range("D379").NumberFormat
$ws.Range("D381").NumberFormat = $ws.Range("D379").NumberFormat
